$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from H1 into the new I1/J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(7,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(5,6),
    @(11,11),
    @(9,9),
    @(8,8),
    @(11,11),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(7,7),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(7,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(6,6),
    @(7,7),
    @(8,8),
    @(7,8),
    @(5,5),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,7),
    @(8,9),
    @(8,9),
    @(8,8),
    @(8,9),
    @(9,9),
    @(5,5),
    @(5,5),
    @(9,9),
    @(5,5),
    @(7,7),
    @(7,7),
    @(6,6)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $r = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}

Write-Host "I0 and IF columns added"
